$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 08:05"

# Kazajistan overtakes Barein in total cases, so rows 59/60 swap labels.
# Row 59 becomes Kazajistan with the newly updated figures.
$ws.Range("A59").Value = "Kazajistan"
$ws.Range("B59").Value = 5240
$ws.Range("C59").Value = 33
$ws.Range("D59").Value = 2108
$ws.Range("E59").Value = 3100
$ws.Range("F59").Value = 33
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 32

# Row 60 becomes Barein, carrying over the previous (unchanged) Barein figures.
$ws.Range("A60").Value = "Barein"
$ws.Range("B60").Value = 5236
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 2152
$ws.Range("E60").Value = 3076
$ws.Range("F60").Value = 6
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 8

# El Salvador figures updated, row position unchanged.
$ws.Range("B100").Value = 998
$ws.Range("C100").Value = 40
$ws.Range("D100").Value = 349
$ws.Range("E100").Value = 631
